$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.069.16'
$ws.Range("E2").Value = '  +2.14%  '

$ws.Range("D3").Value = '2.312.73'
$ws.Range("E3").Value = '  +1.70%  '

$ws.Range("E4").Value = '  +0.10%  '

$c = $ws.Range("D5")
$c.Value = "'302.03"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '

$c = $ws.Range("D6")
$c.Value = "'101.60"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +6.28%  '

$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("E8").Value = '  +0.06%  '

$c = $ws.Range("D9")
$c.Value = "'0.514"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +4.86%  '

$c = $ws.Range("D10")
$c.Value = "'36.22"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +9.51%  '

$c = $ws.Range("D11")
$c.Value = "'0.0792"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '

$ws.Range("E12").Value = '  +2.69%  '

$c = $ws.Range("D13")
$c.Value = "'17.85"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +11.12%  '

$c = $ws.Range("D14")
$c.Value = "'6.87"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.95%  '

$ws.Range("D15").Value = '2.673.62'
$ws.Range("E15").Value = '  +1.85%  '

$ws.Range("D16").Value = '2.324.46'
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("E17").Value = '  +2.51%  '

$ws.Range("D18").Value = '42.965.57'
$ws.Range("E18").Value = '  +2.00%  '

$c = $ws.Range("D19")
$c.Value = "'12.52"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +6.92%  '

$c = $ws.Range("D20")
$c.Value = "'6.22"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.72%  '

$ws.Range("D21").Value = '0.0₃0902'
$ws.Range("E21").Value = '  +1.21%  '

$c = $ws.Range("D22")
$c.Value = "'67.86"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.35%  '

$c = $ws.Range("D23")
$c.Value = "'236.38"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("E24").Value = '  +13.77%  '

$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("E26").Value = '  -0.12%  '

$c = $ws.Range("D27")
$c.Value = "'24.69"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.72%  '

$ws.Range("E28").Value = '  +9.20%  '

$c = $ws.Range("D29")
$c.Value = "'34.81"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.29%  '

$c = $ws.Range("D30")
$c.Value = "'168.71"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.65%  '

$c = $ws.Range("D31")
$c.Value = "'9.17"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.30%  '

$ws.Range("E32").Value = '  +0.02%  '

$c = $ws.Range("D33")
$c.Value = "'5.01"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.25%  '

$c = $ws.Range("D34")
$c.Value = "'4.68"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("E35").Value = '  +3.87%  '

$ws.Range("E36").Value = '  +3.34%  '

$c = $ws.Range("D37")
$c.Value = "'0.0693"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.58%  '

$c = $ws.Range("D38")
$c.Value = "'0.102"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.08%  '

$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("E40").Value = '  +3.86%  '

$ws.Range("E41").Value = '  +0.54%  '

$ws.Range("D42").Value = '1.984.21'
$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("E43").Value = '  +4.27%  '

$c = $ws.Range("D44")
$c.Value = "'2.21"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.29%  '

$c = $ws.Range("D45")
$c.Value = "'10.21"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +6.64%  '

$c = $ws.Range("D46")
$c.Value = "'2.91"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.86%  '

$c = $ws.Range("D47")
$c.Value = "'17.68"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.47%  '

$c = $ws.Range("D48")
$c.Value = "'56.09"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +7.20%  '

$ws.Range("D49").Value = '2.539.79'
$ws.Range("E49").Value = '  +1.76%  '

$c = $ws.Range("D50")
$c.Value = "'1.54"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +4.14%  '

$c = $ws.Range("D51")
$c.Value = "'4.55"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.65%  '
